# Datentreiber.xlsx maintenance edit:
# - Remove the "Testautomatisierer" / "Testmanager" rows of test data,
#   leaving only the styled-but-empty A3/A4 cells behind.
# - Update the current selection on the "Testdaten" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Testdaten")

# Clear the now-obsolete B3/B4 values and the A3/A4 shared-string values
# (the cell formatting/style on column A must stay behind).
$ws.Range("B3").ClearContents()
$ws.Range("B4").ClearContents()
$ws.Range("A3").ClearContents()
$ws.Range("A4").ClearContents()

# Update the active selection to match the new editing state.
$ws.Activate()
$ws.Range("A3:XFD4").Select()
